$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has data rows through row 244 (dates up to 2021-05-02).
# Append three more daily rows (245-247), reusing the formatting (date style
# with border/bold/center alignment, numFmt "YYYY-MM-DD HH:MM:SS") that the
# existing rows already use in column A, by copying the last row's format
# down before writing the new values.
$ws.Range("A244:D244").Copy()
$ws.Range("A245:D247").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(245, 1).Value = 44319
$ws.Cells.Item(245, 2).Value = 0
$ws.Cells.Item(245, 3).Value = 4
$ws.Cells.Item(245, 4).Value = 107.0663811563169

$ws.Cells.Item(246, 1).Value = 44320
$ws.Cells.Item(246, 2).Value = 0
$ws.Cells.Item(246, 3).Value = 3
$ws.Cells.Item(246, 4).Value = 80.29978586723769

$ws.Cells.Item(247, 1).Value = 44321
$ws.Cells.Item(247, 2).Value = 0
$ws.Cells.Item(247, 3).Value = 2
$ws.Cells.Item(247, 4).Value = 53.53319057815846
